$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new year of data (2020) has been added to the table as column N.
# Fill in the values first ...
$ws.Range("N3").Value = 2020
$ws.Range("N4").Value = 15
$ws.Range("N5").Value = 1308.3

# ... then copy the formatting from the previous year's column (M) so the
# new column matches the look (number formats, borders, alignment, etc.)
# of the rest of the table.
$ws.Range("M3:M5").Copy()
$ws.Range("N3:N5").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Reflect the new selection left after entering the data.
$ws.Range("N6").Select()
